# Padronizacao dos parametros de parada do gradiente.
# Updates numeric result cells on rows 2 (secao_aurea) and 3 (armijo).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("secao_aurea" / Rosenbrock) updated results
$ws.Range("D2").Value = [double]"1.944159978522072e-15"
$ws.Range("E2").Value = 442.0
$ws.Range("H2").Value = [double]"5.6592057119481806e-6"
$ws.Range("I2").Value = 5000.0
$ws.Range("K2").Value = 0.5
$ws.Range("L2").Value = [double]"1.501992759363086e-13"
$ws.Range("M2").Value = 2869.0

# Row 3 ("armijo" / Gaussian) updated parameters and results
$ws.Range("B3").Value = 0.5
$ws.Range("C3").Value = 0.001
$ws.Range("D3").Value = [double]"1.6156315450122798e-10"
$ws.Range("E3").Value = 1.0
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.5
$ws.Range("H3").Value = [double]"5.960322351693935e-7"
$ws.Range("I3").Value = 250.0
$ws.Range("J3").Value = 0.3
$ws.Range("K3").Value = 0.9
$ws.Range("L3").Value = [double]"5.957194087948522e-7"
$ws.Range("M3").Value = 278.0
